$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-09 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-10 Saturday", 2) | Out-Null
$d.Content.Find.Execute("83-9=74", $true, $false, $false, $false, $false, $true, 1, $false, "14-13=1", 2) | Out-Null
$d.Content.Find.Execute("80+16=96", $true, $false, $false, $false, $false, $true, 1, $false, "36+23=59", 2) | Out-Null
$d.Content.Find.Execute("74-23=51", $true, $false, $false, $false, $false, $true, 1, $false, "83-74=9", 2) | Out-Null
$d.Content.Find.Execute("91-61=30", $true, $false, $false, $false, $false, $true, 1, $false, "61-33=28", 2) | Out-Null
$d.Content.Find.Execute("16+20=36", $true, $false, $false, $false, $false, $true, 1, $false, "62+11=73", 2) | Out-Null
$d.Content.Find.Execute("11+3=14", $true, $false, $false, $false, $false, $true, 1, $false, "6+78=84", 2) | Out-Null
$d.Content.Find.Execute("12+23=35", $true, $false, $false, $false, $false, $true, 1, $false, "43+51=94", 2) | Out-Null
$d.Content.Find.Execute("28+1=29", $true, $false, $false, $false, $false, $true, 1, $false, "26+39=65", 2) | Out-Null
$d.Content.Find.Execute("24+62=86", $true, $false, $false, $false, $false, $true, 1, $false, "36+41=77", 2) | Out-Null
$d.Content.Find.Execute("31+58=89", $true, $false, $false, $false, $false, $true, 1, $false, "41-12=29", 2) | Out-Null
$d.Content.Find.Execute("56-49=7", $true, $false, $false, $false, $false, $true, 1, $false, "29+16=45", 2) | Out-Null
$d.Content.Find.Execute("71+23=94", $true, $false, $false, $false, $false, $true, 1, $false, "81-1=80", 2) | Out-Null
$d.Content.Find.Execute("67-14=53", $true, $false, $false, $false, $false, $true, 1, $false, "80-29=51", 2) | Out-Null
$d.Content.Find.Execute("53-25=28", $true, $false, $false, $false, $false, $true, 1, $false, "14+35=49", 2) | Out-Null
$d.Content.Find.Execute("99-11=88", $true, $false, $false, $false, $false, $true, 1, $false, "13+15=28", 2) | Out-Null
$d.Content.Find.Execute("20+34=54", $true, $false, $false, $false, $false, $true, 1, $false, "46+35=81", 2) | Out-Null
$d.Content.Find.Execute("8+67=75", $true, $false, $false, $false, $false, $true, 1, $false, "47-20=27", 2) | Out-Null
$d.Content.Find.Execute("84-83=1", $true, $false, $false, $false, $false, $true, 1, $false, "63-38=25", 2) | Out-Null
$d.Content.Find.Execute("47+7=54", $true, $false, $false, $false, $false, $true, 1, $false, "5+9=14", 2) | Out-Null
$d.Content.Find.Execute("73-45=28", $true, $false, $false, $false, $false, $true, 1, $false, "65+33=98", 2) | Out-Null
$d.Content.Find.Execute("47-32=15", $true, $false, $false, $false, $false, $true, 1, $false, "48-27=21", 2) | Out-Null
$d.Content.Find.Execute("81-78=3", $true, $false, $false, $false, $false, $true, 1, $false, "25-5=20", 2) | Out-Null
$d.Content.Find.Execute("81-74=7", $true, $false, $false, $false, $false, $true, 1, $false, "8+20=28", 2) | Out-Null
$d.Content.Find.Execute("4+9=13", $true, $false, $false, $false, $false, $true, 1, $false, "88-40=48", 2) | Out-Null
$d.Content.Find.Execute("19+26=45", $true, $false, $false, $false, $false, $true, 1, $false, "88-46=42", 2) | Out-Null
$d.Content.Find.Execute("13+53=66", $true, $false, $false, $false, $false, $true, 1, $false, "5+50=55", 2) | Out-Null
$d.Content.Find.Execute("87-37=50", $true, $false, $false, $false, $false, $true, 1, $false, "99-23=76", 2) | Out-Null
$d.Content.Find.Execute("29+4=33", $true, $false, $false, $false, $false, $true, 1, $false, "4+87=91", 2) | Out-Null
$d.Content.Find.Execute("36-19=17", $true, $false, $false, $false, $false, $true, 1, $false, "47+38=85", 2) | Out-Null
$d.Content.Find.Execute("58-54=4", $true, $false, $false, $false, $false, $true, 1, $false, "12+16=28", 2) | Out-Null
$d.Content.Find.Execute("83-19=64", $true, $false, $false, $false, $false, $true, 1, $false, "19+54=73", 2) | Out-Null
$d.Content.Find.Execute("89-2=87", $true, $false, $false, $false, $false, $true, 1, $false, "24+61=85", 2) | Out-Null
$d.Content.Find.Execute("73-41=32", $true, $false, $false, $false, $false, $true, 1, $false, "59+17=76", 2) | Out-Null
$d.Content.Find.Execute("52-50=2", $true, $false, $false, $false, $false, $true, 1, $false, "33-0=33", 2) | Out-Null
$d.Content.Find.Execute("40+59=99", $true, $false, $false, $false, $false, $true, 1, $false, "47+25=72", 2) | Out-Null
$d.Content.Find.Execute("54+2=56", $true, $false, $false, $false, $false, $true, 1, $false, "79-79=0", 2) | Out-Null
$d.Content.Find.Execute("21+0=21", $true, $false, $false, $false, $false, $true, 1, $false, "5+23=28", 2) | Out-Null
$d.Content.Find.Execute("8+51=59", $true, $false, $false, $false, $false, $true, 1, $false, "15-11=4", 2) | Out-Null
$d.Content.Find.Execute("61-50=11", $true, $false, $false, $false, $false, $true, 1, $false, "72+1=73", 2) | Out-Null
$d.Content.Find.Execute("32-8=24", $true, $false, $false, $false, $false, $true, 1, $false, "17+14=31", 2) | Out-Null
$d.Content.Find.Execute("18+48=66", $true, $false, $false, $false, $false, $true, 1, $false, "49+44=93", 2) | Out-Null
$d.Content.Find.Execute("28+54=82", $true, $false, $false, $false, $false, $true, 1, $false, "61+5=66", 2) | Out-Null
$d.Content.Find.Execute("98-14=84", $true, $false, $false, $false, $false, $true, 1, $false, "85-38=47", 2) | Out-Null
$d.Content.Find.Execute("88-52=36", $true, $false, $false, $false, $false, $true, 1, $false, "71-53=18", 2) | Out-Null
$d.Content.Find.Execute("36+36=72", $true, $false, $false, $false, $false, $true, 1, $false, "27-12=15", 2) | Out-Null
$d.Content.Find.Execute("17+0=17", $true, $false, $false, $false, $false, $true, 1, $false, "89+4=93", 2) | Out-Null
$d.Content.Find.Execute("90-38=52", $true, $false, $false, $false, $false, $true, 1, $false, "66-30=36", 2) | Out-Null
$d.Content.Find.Execute("16+14=30", $true, $false, $false, $false, $false, $true, 1, $false, "67-38=29", 2) | Out-Null
$d.Content.Find.Execute("36+49=85", $true, $false, $false, $false, $false, $true, 1, $false, "53+8=61", 2) | Out-Null
$d.Content.Find.Execute("1+91=92", $true, $false, $false, $false, $false, $true, 1, $false, "47-34=13", 2) | Out-Null
$d.Content.Find.Execute("92-12=80", $true, $false, $false, $false, $false, $true, 1, $false, "33-22=11", 2) | Out-Null
$d.Content.Find.Execute("34+48=82", $true, $false, $false, $false, $false, $true, 1, $false, "80-4=76", 2) | Out-Null
$d.Content.Find.Execute("11+22=33", $true, $false, $false, $false, $false, $true, 1, $false, "57+8=65", 2) | Out-Null
$d.Content.Find.Execute("79-76=3", $true, $false, $false, $false, $false, $true, 1, $false, "67-43=24", 2) | Out-Null
$d.Content.Find.Execute("0+84=84", $true, $false, $false, $false, $false, $true, 1, $false, "31-2=29", 2) | Out-Null
$d.Content.Find.Execute("19+46=65", $true, $false, $false, $false, $false, $true, 1, $false, "98-72=26", 2) | Out-Null
$d.Content.Find.Execute("21+20=41", $true, $false, $false, $false, $false, $true, 1, $false, "84-3=81", 2) | Out-Null
$d.Content.Find.Execute("80-45=35", $true, $false, $false, $false, $false, $true, 1, $false, "56-35=21", 2) | Out-Null
$d.Content.Find.Execute("65-18=47", $true, $false, $false, $false, $false, $true, 1, $false, "29+43=72", 2) | Out-Null
$d.Content.Find.Execute("38-19=19", $true, $false, $false, $false, $false, $true, 1, $false, "81-38=43", 2) | Out-Null
$d.Content.Find.Execute("82-66=16", $true, $false, $false, $false, $false, $true, 1, $false, "69+10=79", 2) | Out-Null
$d.Content.Find.Execute("84-20=64", $true, $false, $false, $false, $false, $true, 1, $false, "47-17=30", 2) | Out-Null
$d.Content.Find.Execute("1+40=41", $true, $false, $false, $false, $false, $true, 1, $false, "19+62=81", 2) | Out-Null
$d.Content.Find.Execute("42+9=51", $true, $false, $false, $false, $false, $true, 1, $false, "48+41=89", 2) | Out-Null
$d.Content.Find.Execute("20+58=78", $true, $false, $false, $false, $false, $true, 1, $false, "30-3=27", 2) | Out-Null
$d.Content.Find.Execute("21+19=40", $true, $false, $false, $false, $false, $true, 1, $false, "40+11=51", 2) | Out-Null
$d.Content.Find.Execute("62+0=62", $true, $false, $false, $false, $false, $true, 1, $false, "35-5=30", 2) | Out-Null
$d.Content.Find.Execute("93-85=8", $true, $false, $false, $false, $false, $true, 1, $false, "15-12=3", 2) | Out-Null
$d.Content.Find.Execute("10+64=74", $true, $false, $false, $false, $false, $true, 1, $false, "2+11=13", 2) | Out-Null
$d.Content.Find.Execute("77-41=36", $true, $false, $false, $false, $false, $true, 1, $false, "20+67=87", 2) | Out-Null
$d.Content.Find.Execute("38+57=95", $true, $false, $false, $false, $false, $true, 1, $false, "76-9=67", 2) | Out-Null
$d.Content.Find.Execute("45+21=66", $true, $false, $false, $false, $false, $true, 1, $false, "38-29=9", 2) | Out-Null
$d.Content.Find.Execute("43+1=44", $true, $false, $false, $false, $false, $true, 1, $false, "44-11=33", 2) | Out-Null
$d.Content.Find.Execute("85-52=33", $true, $false, $false, $false, $false, $true, 1, $false, "48+12=60", 2) | Out-Null
$d.Content.Find.Execute("74-1=73", $true, $false, $false, $false, $false, $true, 1, $false, "33+46=79", 2) | Out-Null
$d.Content.Find.Execute("67+24=91", $true, $false, $false, $false, $false, $true, 1, $false, "96-25=71", 2) | Out-Null
$d.Content.Find.Execute("16-5=11", $true, $false, $false, $false, $false, $true, 1, $false, "87-49=38", 2) | Out-Null
$d.Content.Find.Execute("2+20=22", $true, $false, $false, $false, $false, $true, 1, $false, "59-46=13", 2) | Out-Null
$d.Content.Find.Execute("93-44=49", $true, $false, $false, $false, $false, $true, 1, $false, "56+25=81", 2) | Out-Null
$d.Content.Find.Execute("34-10=24", $true, $false, $false, $false, $false, $true, 1, $false, "3+89=92", 2) | Out-Null
$d.Content.Find.Execute("2+15=17", $true, $false, $false, $false, $false, $true, 1, $false, "12-5=7", 2) | Out-Null
$d.Content.Find.Execute("78-22=56", $true, $false, $false, $false, $false, $true, 1, $false, "15+32=47", 2) | Out-Null
$d.Content.Find.Execute("7+83=90", $true, $false, $false, $false, $false, $true, 1, $false, "48+51=99", 2) | Out-Null
$d.Content.Find.Execute("66-16=50", $true, $false, $false, $false, $false, $true, 1, $false, "21+37=58", 2) | Out-Null
$d.Content.Find.Execute("85-39=46", $true, $false, $false, $false, $false, $true, 1, $false, "91-41=50", 2) | Out-Null
$d.Content.Find.Execute("93-34=59", $true, $false, $false, $false, $false, $true, 1, $false, "31+60=91", 2) | Out-Null
$d.Content.Find.Execute("56-31=25", $true, $false, $false, $false, $false, $true, 1, $false, "55-11=44", 2) | Out-Null
$d.Content.Find.Execute("89-66=23", $true, $false, $false, $false, $false, $true, 1, $false, "71-31=40", 2) | Out-Null
$d.Content.Find.Execute("36+21=57", $true, $false, $false, $false, $false, $true, 1, $false, "77+13=90", 2) | Out-Null
$d.Content.Find.Execute("89-7=82", $true, $false, $false, $false, $false, $true, 1, $false, "29-19=10", 2) | Out-Null
$d.Content.Find.Execute("8+5=13", $true, $false, $false, $false, $false, $true, 1, $false, "31-2=29", 2) | Out-Null
$d.Content.Find.Execute("2+70=72", $true, $false, $false, $false, $false, $true, 1, $false, "45-11=34", 2) | Out-Null
$d.Content.Find.Execute("67+14=81", $true, $false, $false, $false, $false, $true, 1, $false, "29+54=83", 2) | Out-Null
$d.Content.Find.Execute("22+5=27", $true, $false, $false, $false, $false, $true, 1, $false, "63-27=36", 2) | Out-Null
$d.Content.Find.Execute("94-35=59", $true, $false, $false, $false, $false, $true, 1, $false, "88-60=28", 2) | Out-Null
$d.Content.Find.Execute("9+44=53", $true, $false, $false, $false, $false, $true, 1, $false, "7+8=15", 2) | Out-Null
$d.Content.Find.Execute("76-73=3", $true, $false, $false, $false, $false, $true, 1, $false, "80-55=25", 2) | Out-Null
$d.Content.Find.Execute("87+3=90", $true, $false, $false, $false, $false, $true, 1, $false, "10+17=27", 2) | Out-Null
$d.Content.Find.Execute("14+12=26", $true, $false, $false, $false, $false, $true, 1, $false, "20+56=76", 2) | Out-Null
$d.Content.Find.Execute("38-27=11", $true, $false, $false, $false, $false, $true, 1, $false, "21+17=38", 2) | Out-Null
